# Append: 2025-10-28 12:36 JST
# Updates the "ランサーズ" (Lancers) sheet:
#   - A2:A15  retrieval timestamp 06:37:23 -> 12:36:23
#   - B15     title text for the new/replaced listing
#   - D15     price range text for the new/replaced listing
#   - F15     URL text + underlying hyperlink target for the new/replaced listing
#
# NOTE on hyperlinks: this host's Hyperlinks.Delete() (called on *any* scoped
# Range.Hyperlinks collection, even a single cell) clears every hyperlink
# relationship on the sheet once the workbook is saved, not just the scoped
# one - unless every surviving link is explicitly re-added with
# Hyperlinks.Add() before save. So instead of deleting just the F15 link (which
# would silently drop the F2:F14 links too), every URL cell's hyperlink is
# rebuilt from its current address, with only F15's address replaced.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-10-28 06:37:23"
$newTimestamp = "2025-10-28 12:36:23"

# --- A2:A15 timestamp refresh -------------------------------------------------
for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}

# --- Row 15: listing was replaced by a new one --------------------------------
$ws.Range("B15").Value = "AWS構成(EC2・RDS・WAF・CloudFront)トラブル調査に向けた概要提案募集"
$ws.Range("D15").Value = "~ 5,000 円 / 固定"

$newUrl = "https://www.lancers.jp/work/detail/5422200"

# --- Rebuild the F2:F15 hyperlinks, changing only F15's target ---------------
$lastRow = 15
$firstRow = 2

# Capture the existing "Hyperlink" cell style once, up front, so it can be
# restored after Hyperlinks.Add() re-applies its own (duplicate) style record.
$hyperlinkStyle = $ws.Range("F2").Style

$urls = @{}
for ($row = $firstRow; $row -le $lastRow; $row++) {
    $urls[$row] = $ws.Cells.Item($row, 6).Value2
}
$urls[15] = $newUrl

$ws.Hyperlinks.Delete()

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value = $urls[$row]
    $ws.Hyperlinks.Add($cell, $urls[$row])
    $cell.Style = $hyperlinkStyle
}
